$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555999625168"
$ws1.Range("B2").Value = "go_stims-16512555999305174.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555999435182.csv"
$ws1.Range("B4").Value = "go_stims-16512555999445162.csv"
$ws1.Range("B5").Value = "GNG_stims-165125559996052.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651255601335645"
$ws2.Range("B2").Value = "ZB-match_8-16512555999785135.csv"
$ws2.Range("B3").Value = "OB-16512556009546444.csv"
$ws2.Range("B4").Value = "ZB-match_3-1651255600128664.csv"
$ws2.Range("B5").Value = "TB-1651255601022649.csv"
$ws2.Range("B6").Value = "OB-16512556004236436.csv"
$ws2.Range("B7").Value = "ZB-match_0-16512556002586462.csv"
$ws2.Range("B8").Value = "TB-1651255601093645.csv"
$ws2.Range("B9").Value = "OB-16512556003676436.csv"
$ws2.Range("B10").Value = "TB-16512556013146465.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512556013376544"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512556013846514"
$ws4.Range("B2").Value = "MM_stims-16512556013516526.csv"
$ws4.Range("B3").Value = "ZM_stims-16512556013396444.csv"
$ws4.Range("B4").Value = "MM_stims-16512556013676443.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556013526456.csv"
$ws4.Range("B6").Value = "MM_stims-16512556013836467.csv"
$ws4.Range("B7").Value = "ZM_stims-16512556013686452.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512556014636457"
$ws5.Range("B2").Value = "SAT_stims-16512556013896458.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512556014476454.csv"
$ws5.Range("B4").Value = "SAT_stims-1651255601415647.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512556014326515.csv"
